# Apply "Minor updates to Excel workbooks" edit to the decision matrix.
# All of the affected data lives on the "solution" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solution")

# --- Raw input values that changed (everything else on the sheet is
#     formula-driven and will recompute automatically) ---

# Row 2: weight inputs for cost / performance / features
$ws.Range("B2").Value = 33
$ws.Range("D2").Value = 33
$ws.Range("F2").Value = 100

# Row 5: weight inputs for the six sub-criteria
$ws.Range("B5").Value = 100
$ws.Range("C5").Value = 67
$ws.Range("D5").Value = 43
$ws.Range("E5").Value = 100
$ws.Range("F5").Value = 100
$ws.Range("G5").Value = 100

# --- Text corrections / casing fixes ---

# "comform" -> "comfort"
$ws.Range("B4").Value = "comfort"

# "Scores" -> "scores"
$ws.Range("B9").Value = "scores"

# "Rating" -> "rating"
$ws.Range("H9").Value = "rating"

# "Rank" -> "rank"
$ws.Range("I9").Value = "rank"

# "Legend" -> "legend"
$ws.Range("K2").Value = "legend"

# Remove the stray "Alternative" label above the second table
$ws.Range("A9").ClearContents()
